# Update Case_2_160 vm_pu.xlsx results: bus voltage magnitudes for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Per-row new values for columns B,C,D,E,F,I,J,K,L,M,N (G/H are unchanged)
$newValues = @(
    @{ Row = 2; Values = @(1.02,1.037340378303135,1.044178949661174,1.05221432827621,1.057584219475618,1.036712514000534,1.042443781575285,1.046950728911875,1.054963680062444,1.060318789219502,1.01803950637945) },
    @{ Row = 3; Values = @(1.02,1.038511074985778,1.045074024077472,1.053280907062049,1.058674644738434,1.036960149036202,1.043257705240028,1.047656769698613,1.055842408349053,1.061222377196264,1.018312679113075) },
    @{ Row = 4; Values = @(1.02,1.039268196578572,1.045652592455614,1.053971610629533,1.059380348955775,1.037118381864645,1.043783453724002,1.048112389344859,1.05641096069947,1.061806590366119,1.018489037336832) },
    @{ Row = 5; Values = @(1.02,1.03958639678062,1.045895678054118,1.054262115736919,1.059677057626502,1.037184423262094,1.044004260037325,1.048303636242602,1.056649969949004,1.062052081746473,1.018563082036327) },
    @{ Row = 6; Values = @(1.02,1.039639818540728,1.045936484725023,1.054310900649037,1.059726878125298,1.037195483773789,1.044041321602021,1.048335730112266,1.056690100066101,1.062093294281801,1.018575508824081) },
    @{ Row = 7; Values = @(1.02,1.039272448751704,1.045655841145027,1.053975491848777,1.059384313470746,1.037119266197899,1.043786405006007,1.048114945956444,1.056414154392337,1.061809871072962,1.018490027103955) },
    @{ Row = 8; Values = @(1.02,1.037736103886031,1.044481569565338,1.05257466866107,1.05795270792921,1.036796618270019,1.042719040614215,1.047189594790561,1.055260659602133,1.060624258270588,1.018131909867804) },
    @{ Row = 9; Values = @(1.02,1.035025759178714,1.042407712327353,1.050110482731557,1.05543098277966,1.036212726342805,1.040831176303894,1.045549528037066,1.053227707083677,1.058531450371361,1.017497770284943) },
    @{ Row = 10; Values = @(1.02,1.033216678861808,1.041022003014896,1.048470523655201,1.053750438357903,1.035813145904611,1.039567830804744,1.044449754002666,1.051872150314914,1.057133792859839,1.017072922285086) },
    @{ Row = 11; Values = @(1.02,1.032432787335734,1.040421226257028,1.047761069058451,1.053022879753895,1.035637674467539,1.039019646884488,1.04397201660438,1.051285113050157,1.0565280029831,1.016888460550934) },
    @{ Row = 12; Values = @(1.02,1.032141530824292,1.040197956951801,1.047497644356891,1.05275265054933,1.035572128218073,1.038815853659719,1.043794333464679,1.051067049824915,1.056302895536532,1.016819867906678) },
    @{ Row = 13; Values = @(1.02,1.032204010181367,1.04024585414668,1.047554145410065,1.052810614808946,1.03558620478258,1.038859575856774,1.043832457528606,1.051113825654586,1.05635118593759,1.01683458467538) },
    @{ Row = 14; Values = @(1.02,1.032408713712812,1.040402773076866,1.047739292288205,1.053000542151765,1.035632263909354,1.039002804814148,1.043957333960325,1.051267088112913,1.056509397373878,1.016882792200062) },
    @{ Row = 15; Values = @(1.02,1.032534827084798,1.040499440876548,1.047853380401029,1.053117565077283,1.035660593642453,1.039091029907131,1.044034243928581,1.051361516637035,1.056606864698548,1.016912484455505) },
    @{ Row = 16; Values = @(1.02,1.033268691449017,1.041061858624591,1.048517621642387,1.053798726658973,1.035824739721381,1.039604187716425,1.044481427613438,1.051911108472251,1.05717398454188,1.017085153851509) },
    @{ Row = 17; Values = @(1.02,1.033728877106606,1.041414445759183,1.048934458496244,1.054226035082974,1.035927048027658,1.039925770277971,1.044761524593153,1.05225583315203,1.057529564038731,1.017193330821349) },
    @{ Row = 18; Values = @(1.02,1.033997242895603,1.041620031016606,1.049177655959813,1.054475289285098,1.035986486310752,1.040113233346542,1.044924753038896,1.052456898297129,1.057736910190101,1.017256380417593) },
    @{ Row = 19; Values = @(1.02,1.034088739811994,1.0416901179359,1.049260590812017,1.054560280658475,1.036006713133526,1.040177134678272,1.044980384740541,1.052525455157886,1.057807600162751,1.017277870532714) },
    @{ Row = 20; Values = @(1.02,1.033679509009266,1.041376624036919,1.048889729275723,1.054180187654393,1.035916095770526,1.039891278960263,1.044731488087959,1.052218848178027,1.057491419662308,1.01718172944399) },
    @{ Row = 21; Values = @(1.02,1.032348435967935,1.040356567516449,1.04768476843795,1.052944612780189,1.035618710814187,1.038960632248238,1.043920567323415,1.051221956459704,1.056462810542677,1.016868598367537) },
    @{ Row = 22; Values = @(1.02,1.031511048813188,1.039714557247294,1.046927730301591,1.052167864188438,1.035429601860523,1.038374494515878,1.043409377072474,1.050595104540851,1.055815561427657,1.016671284388371) },
    @{ Row = 23; Values = @(1.02,1.031955010579558,1.040054961813973,1.047328996852831,1.052579623503016,1.035530054175829,1.038685312642675,1.043680495187105,1.050927417109806,1.056158730130177,1.016775925682117) },
    @{ Row = 24; Values = @(1.02,1.033701816490389,1.041393714270464,1.048909940288871,1.054200904095667,1.035921045354669,1.039906864443909,1.044745060747264,1.05223556011734,1.05750865564147,1.017186971756019) },
    @{ Row = 25; Values = @(1.02,1.035726826270536,1.042944405855067,1.050747033177167,1.056082800126429,1.036365495047791,1.041320072027453,1.045974650287754,1.053753317395288,1.059072921566171,1.017662077904109) }
)

$colLetters = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

foreach ($entry in $newValues) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $cellRef = "$($colLetters[$i])$r"
        $ws.Range($cellRef).Value = $vals[$i]
    }
}
